# Apply updated cryptocurrency price/volume data (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$q = "'"   # leading apostrophe forces Excel to store numeric-looking text as text

$ws.Range("D2").Value = '44.017.11'
$ws.Range("E2").Value = '  +1.76%  '

$ws.Range("D3").Value = '2.375.29'
$ws.Range("E3").Value = '  +0.73%  '

$ws.Range("D4").Value = $q + '1.01'
$ws.Range("E4").Value = '  +0.49%  '

$ws.Range("D5").Value = $q + '0.687'
$ws.Range("E5").Value = '  +5.53%  '

$ws.Range("D6").Value = $q + '241.64'
$ws.Range("E6").Value = '  +3.44%  '

$ws.Range("D7").Value = $q + '75.82'
$ws.Range("E7").Value = '  +6.00%  '

$ws.Range("E8").Value = '  +0.14%  '

$ws.Range("D9").Value = $q + '0.625'
$ws.Range("E9").Value = '  +29.09%  '

$ws.Range("D10").Value = $q + '0.102'
$ws.Range("E10").Value = '  +4.19%  '

$ws.Range("D11").Value = $q + '57.29'
$ws.Range("E11").Value = '  +0.87%  '

$ws.Range("D12").Value = $q + '32.88'
$ws.Range("E12").Value = '  +21.13%  '

$ws.Range("D13").Value = $q + '7.49'
$ws.Range("E13").Value = '  +18.88%  '

$ws.Range("D14").Value = $q + '0.108'
$ws.Range("E14").Value = '  +1.43%  '

$ws.Range("D15").Value = '2.730.37'
$ws.Range("E15").Value = '  +0.73%  '

$ws.Range("D16").Value = $q + '16.92'
$ws.Range("E16").Value = '  +4.79%  '

$ws.Range("D17").Value = $q + '0.920'
$ws.Range("E17").Value = '  +6.30%  '

$ws.Range("D18").Value = '2.375.42'
$ws.Range("E18").Value = '  +1.36%  '

$ws.Range("D19").Value = '44.079.83'
$ws.Range("E19").Value = '  +1.87%  '

$ws.Range("D20").Value = $q + '0.0000102'
$ws.Range("E20").Value = '  +0.89%  '

$ws.Range("D21").Value = $q + '6.65'
$ws.Range("E21").Value = '  +4.91%  '

$ws.Range("D22").Value = $q + '77.92'
$ws.Range("E22").Value = '  +4.70%  '

$ws.Range("D23").Value = $q + '256.12'
$ws.Range("E23").Value = '  +2.43%  '

$ws.Range("E24").Value = '  -0.02%  '

$ws.Range("D25").Value = $q + '2.53'
$ws.Range("E25").Value = '  +3.32%  '

$ws.Range("B26").Value = 'WEMIXToken'
$ws.Range("C26").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D26").Value = $q + '3.66'
$ws.Range("E26").Value = '  -3.96%  '

$ws.Range("B27").Value = 'Cosmos'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D27").Value = $q + '11.03'
$ws.Range("E27").Value = '  +9.99%  '

$ws.Range("D28").Value = $q + '1.76'
$ws.Range("E28").Value = '  +16.69%  '

$ws.Range("E29").Value = '  +5.44%  '

$ws.Range("D30").Value = $q + '23.16'
$ws.Range("E30").Value = '  +2.76%  '

$ws.Range("D31").Value = $q + '174.98'
$ws.Range("E31").Value = '  +1.18%  '

$ws.Range("D32").Value = $q + '0.129'
$ws.Range("E32").Value = '  -0.24%  '

$ws.Range("D33").Value = $q + '0.135'
$ws.Range("E33").Value = '  +5.47%  '

$ws.Range("D34").Value = $q + '5.31'
$ws.Range("E34").Value = '  +6.79%  '

$ws.Range("D35").Value = $q + '0.0753'
$ws.Range("E35").Value = '  +9.18%  '

$ws.Range("D36").Value = $q + '5.32'
$ws.Range("E36").Value = '  +4.93%  '

$ws.Range("D37").Value = $q + '3.83'
$ws.Range("E37").Value = '  +3.60%  '

$ws.Range("D38").Value = $q + '2.47'
$ws.Range("E38").Value = '  +1.51%  '

$ws.Range("D39").Value = $q + '6.50'
$ws.Range("E39").Value = '  -0.89%  '

$ws.Range("D40").Value = $q + '0.0274'
$ws.Range("E40").Value = '  +7.92%  '

$ws.Range("D41").Value = $q + '8.98'
$ws.Range("E41").Value = '  +0.79%  '

$ws.Range("D42").Value = $q + '18.86'
$ws.Range("E42").Value = '  +1.15%  '

$ws.Range("B43").Value = 'BinanceUSD'
$ws.Range("C43").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D43").Value = $q + '1.00'
$ws.Range("E43").Value = '  +0.16%  '

$ws.Range("B44").Value = 'Algorand'
$ws.Range("C44").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D44").Value = $q + '0.200'
$ws.Range("E44").Value = '  +18.66%  '

$ws.Range("D45").Value = $q + '2.53'
$ws.Range("E45").Value = '  +14.02%  '

$ws.Range("D46").Value = $q + '1.21'
$ws.Range("E46").Value = '  +3.63%  '

$ws.Range("D47").Value = $q + '0.101'
$ws.Range("E47").Value = '  +5.87%  '

$ws.Range("B48").Value = 'TrustWalletToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D48").Value = $q + '1.26'
$ws.Range("E48").Value = '  +4.28%  '

$ws.Range("B49").Value = 'Aave'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D49").Value = $q + '102.28'
$ws.Range("E49").Value = '  +3.15%  '

$ws.Range("D50").Value = $q + '4.51'
$ws.Range("E50").Value = '  +0.56%  '

$ws.Range("D51").Value = $q + '54.65'
$ws.Range("E51").Value = '  +8.85%  '
